# Update latest output (run 11)
$wb = $excel.ActiveWorkbook

# --- Schedule sheet: refresh summary cost figures ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 1441.930581
$schedule.Range("F2").Value = 23.84144479166666

# --- Detailed sheet: refresh price forecast/historical values ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B15").Value = 79.09635
$detailed.Range("B16").Value = 56.98

$detailed.Range("B17").Value = 49.82038
$detailed.Range("C17").Value = "historical"

$detailed.Range("B18").Value = 36.06
$detailed.Range("C18").Value = "historical"

$detailed.Range("B19").Value = 0.08205
$detailed.Range("B20").Value = 24.47133
$detailed.Range("B21").Value = 0.0223
$detailed.Range("B22").Value = -2.98782
$detailed.Range("B23").Value = 0.51
$detailed.Range("B24").Value = 0
$detailed.Range("B25").Value = 0
$detailed.Range("B26").Value = 22.07
$detailed.Range("B27").Value = 4.64073
$detailed.Range("B28").Value = -0.3133
$detailed.Range("B29").Value = -0.31283
$detailed.Range("B30").Value = 52.54048
$detailed.Range("B31").Value = 49.62995
$detailed.Range("B32").Value = 48.42888
$detailed.Range("B33").Value = 36.0601
$detailed.Range("B34").Value = 43.92262
$detailed.Range("B35").Value = 47.63094
$detailed.Range("B36").Value = 57.33246
$detailed.Range("B37").Value = 28.08698
$detailed.Range("B38").Value = 61.4895
$detailed.Range("B39").Value = 70.25961
$detailed.Range("B40").Value = 120.01
$detailed.Range("B41").Value = 126.73
$detailed.Range("B43").Value = 120.01
$detailed.Range("B45").Value = 85.95
$detailed.Range("B46").Value = 71.40000000000001
$detailed.Range("B47").Value = 57.09
$detailed.Range("B48").Value = 58.51148
$detailed.Range("B49").Value = 60.19004
